$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 238, shifting existing rows 238:321 down to 239:322
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row 238 with the new data record
$ws.Cells.Item(238, 1).Value = 5
$ws.Cells.Item(238, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(238, 3).Value = 'Maule'
$ws.Cells.Item(238, 4).Value = 44588
$ws.Cells.Item(238, 5).Value = 7
$ws.Cells.Item(238, 6).Value = 100112043
$ws.Cells.Item(238, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(238, 8).Value = 'Sin especificar'
$ws.Cells.Item(238, 9).Value = 'Primera'
$ws.Cells.Item(238, 10).Value = 400
$ws.Cells.Item(238, 11).Value = 9000
$ws.Cells.Item(238, 12).Value = 9000
$ws.Cells.Item(238, 13).Value = 9000
$ws.Cells.Item(238, 14).Value = '$/caja 80 unidades'
$ws.Cells.Item(238, 15).Value = 'Región del Maule'
$ws.Cells.Item(238, 16).Value = 112
$ws.Cells.Item(238, 17).Value = 80
$ws.Cells.Item(238, 18).Value = 'Hortaliza'
